# Generate Report for Handback
# Update the generated timestamps recorded on the handback status report.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the first file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 19:04:43"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-17 19:04:38"
$wsZhCn.Range("K2").Value = "2016-08-17 19:04:56"

# "de-de" sheet: Correspond Handback DateTime for the first file.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-17 19:05:11"
